$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before I - shifts old column I ("Construcción Iteración 2")
#    into J, leaving a blank column I formatted like its left neighbour (H).
$ws.Columns("I").Insert()

# 2. Re-apply the original column-I formatting (fill/border per row) onto the
#    new column I by copying formats from J (which now holds that formatting).
$ws.Range("J1:J23").Copy()
$ws.Range("I1:I23").PasteSpecial(-4122)

# 3. Header + new "Codigo" values for the Construcción Iteración 3 (C201-C215) items.
$ws.Range("I1").Value = "Codigo"
$codes = @("C201","C202","C203","C204","C205","C206","C207","C208","C209","C210","C211","C212","C213","C214","C215")
for ($i = 0; $i -lt $codes.Length; $i++) {
    $ws.Cells.Item(2 + $i, 9).Value = $codes[$i]
}

# 4. New content label in J2.
$ws.Range("J2").Value = "Informe de Revisión tecnica formal."

# 5. I2 and I16 lose their top border (matches the original table's "section"
#    look - no divider directly under the thick header / above the footer row).
$ws.Range("I2").Borders.Item(8).LineStyle = -4142
$ws.Range("I16").Borders.Item(8).LineStyle = -4142

# 6. Column widths: narrow "Codigo" column, wide new label column.
$ws.Columns("I").ColumnWidth = 7.1666666666666667
$ws.Columns("J").ColumnWidth = 32.6666666666666667

# 7. Update the view: scrolled to show the new columns, selection on K5.
$ws.Application.ActiveWindow.ScrollColumn = 8
$null = $ws.Range("K5").Select()
